# "Generate Report for Handoff" — refresh the localization-status report.
#
# A fresh round of handoff xliff generation completed for the six
# "Ready for handoff" rows (7, 8, 9, 12, 13, 14) on both locale sheets:
#   - the handoff timestamp advances a few seconds, and
#   - the Priority column now records the handoff type ("ht").
# The Overview sheet mirrors the same refreshed "Latest HO Xliff Generate
# Date" for those rows.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 12, 13, 14)

$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-08-30 16:25:55"
}

$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("H$r").Value = "2016-08-30 16:25:51"
    $zhcn.Range("E$r").Value = "ht"
}

$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("H$r").Value = "2016-08-30 16:25:55"
    $dede.Range("E$r").Value = "ht"
}
